# Actualización automática 2025-10-01 08:30:09
#
# Monthly rollover:
#  - "VENTA MENSUAL": month columns (C:junio .. F:septiembre) roll forward
#    one month (C:julio .. F:octubre). The data that lived in D:F shifts
#    left into C:E, and the now-vacated F column (new month, no sales
#    posted yet) is reset to 0. Column widths for the shifted E/F columns
#    follow the data.
#  - "VENTAS POR GRUPO": the figures that had been posted against the
#    (now rolled-off) "septiembre" column for a handful of product
#    groups/clients are cleared back to 0, including the "x de 46" tally
#    row at the bottom of the sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL": shift monthly columns C<-D<-E<-F, zero out F
# ---------------------------------------------------------------------
$ventaMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Header labels move one month forward.
$ventaMensual.Range("C1").Value = "julio"
$ventaMensual.Range("D1").Value = "agosto"
$ventaMensual.Range("E1").Value = "septiembre"
$ventaMensual.Range("F1").Value = "octubre"

# Shift the data block (rows 2-48) left by one column.
$shifted = $ventaMensual.Range("D2:F48").Value2
$ventaMensual.Range("C2:E48").Value2 = $shifted
$ventaMensual.Range("F2:F48").Value2 = 0

# Column widths: column E (new "septiembre") widens to match the old F
# width, column F (new "octubre", freshly emptied) narrows.
$ventaMensual.Columns.Item(5).ColumnWidth = 15.166666666666666
$ventaMensual.Columns.Item(6).ColumnWidth = 12.166666666666666

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": clear the figures tied to the rolled-off month
# ---------------------------------------------------------------------
$ventasPorGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$cellsToClear = @(
    "L9", "M9",
    "D10",
    "L12",
    "D13",
    "L16",
    "L25",
    "L28",
    "L36",
    "I39",
    "L42", "M42",
    "R43",
    "E44", "M44",
    "M45",
    "E47", "L47"
)
foreach ($ref in $cellsToClear) {
    $ventasPorGrupo.Range($ref).Value2 = 0
}

# Bottom tally row ("x de 46") for the columns that were fully cleared.
$ventasPorGrupo.Range("D48").Value = "0 de 46"
$ventasPorGrupo.Range("E48").Value = "0 de 46"
$ventasPorGrupo.Range("I48").Value = "0 de 46"
$ventasPorGrupo.Range("L48").Value = "0 de 46"
$ventasPorGrupo.Range("M48").Value = "0 de 46"
